$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price-like text (e.g. "1.010", "27.648.98") that Excel would
# auto-coerce to a Number when assigned directly, since many look like valid
# floats. Force the whole Price column to Text first, write the new values,
# then clear the temporary formatting so the cell style reverts to the default
# (matching the original un-styled inline-string cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.616.46"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.843.00"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  -2.37%  "
$ws.Range("D5").Value = "316.62"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  -2.03%  "
$ws.Range("E7").Value = "  -2.00%  "
$ws.Range("D8").Value = "0.3733"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").Value = "0.07301"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "0.8707"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").Value = "21.40"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "1.844.93"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "6.701"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "5.395"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").Value = "0.07099"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "88.66"
$ws.Range("E16").Value = "  +4.24%  "
$ws.Range("D17").Value = "1.012"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "0.000008976"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").Value = "15.35"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").Value = "27.636.17"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "5.186"
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").Value = "10.99"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").Value = "2.071.90"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").Value = "1.968"
$ws.Range("E25").Value = "  -3.63%  "
$ws.Range("D26").Value = "154.55"
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("D27").Value = "18.51"
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("D28").Value = "2.152"
$ws.Range("E28").Value = "  +7.76%  "
$ws.Range("D29").Value = "5.317"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "117.39"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "0.08897"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").Value = "1.217"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "0.7723"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").Value = "4.515"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").Value = "2.896"
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("D36").Value = "1.008"
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("D37").Value = "1.127"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").Value = "0.01969"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "0.05297"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "7.141"
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.876"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("D42").Value = "0.1685"
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").Value = "0.5107"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").Value = "8.731"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "10.65"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").Value = "107.01"
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("D47").Value = "0.4735"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").Value = "0.06443"
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("D49").Value = "1.008"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").Value = "1.680"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").Value = "1.839"
$ws.Range("E51").Value = "  -2.64%  "

$ws.Range("D2:D51").ClearFormats()
